$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (rows 4-5, columns E & F) ---
$wsSchedule.Range("E4").Value = 575.3266837499999
$wsSchedule.Range("F4").Value = 30.44056527777778
$wsSchedule.Range("E5").Value = -217.64118675
$wsSchedule.Range("F5").Value = -6.397448170194004

# --- Detailed sheet updates: Price (column B) ---
$wsDetailed.Range("B41").Value = 57.3
$wsDetailed.Range("B42").Value = 59.1531
$wsDetailed.Range("B43").Value = 76.31898
$wsDetailed.Range("B44").Value = 66.49145
$wsDetailed.Range("B45").Value = 75.21257
$wsDetailed.Range("B46").Value = 62.62235
$wsDetailed.Range("B48").Value = 63.03609
$wsDetailed.Range("B49").Value = 57.06003
$wsDetailed.Range("B58").Value = 63.72093
$wsDetailed.Range("B59").Value = 67.19808999999999
$wsDetailed.Range("B60").Value = 67.36196
$wsDetailed.Range("B61").Value = 73.20005
$wsDetailed.Range("B62").Value = 80.61644
$wsDetailed.Range("B63").Value = 63.26957
$wsDetailed.Range("B64").Value = 20.85929
$wsDetailed.Range("B65").Value = 0.05519
$wsDetailed.Range("B66").Value = -2.54451
$wsDetailed.Range("B67").Value = -6
$wsDetailed.Range("B68").Value = -6.89561
$wsDetailed.Range("B69").Value = -7.36813
$wsDetailed.Range("B70").Value = -9.383990000000001
$wsDetailed.Range("B71").Value = -9.699149999999999
$wsDetailed.Range("B73").Value = -15.16775
$wsDetailed.Range("B74").Value = -15.16905
$wsDetailed.Range("B76").Value = -21.6667
$wsDetailed.Range("B78").Value = -21.86572
$wsDetailed.Range("B79").Value = -23.5
$wsDetailed.Range("B80").Value = -24.23063
$wsDetailed.Range("B81").Value = -23.14497
$wsDetailed.Range("B82").Value = -6.88281
$wsDetailed.Range("B83").Value = -5.51
$wsDetailed.Range("B85").Value = 47.13993
$wsDetailed.Range("B87").Value = 63.19697
$wsDetailed.Range("B88").Value = 73.2
$wsDetailed.Range("B89").Value = 79.95
$wsDetailed.Range("B90").Value = 78
$wsDetailed.Range("B91").Value = 73.19
$wsDetailed.Range("B92").Value = 65
$wsDetailed.Range("B94").Value = 59.9297
$wsDetailed.Range("B95").Value = 63.57467

# --- Detailed sheet updates: Type (column C) ---
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("C44").Value = "historical"
